$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tendons(L)")

# Update tendon length values in column B
$ws.Range("B2").Value = 7
$ws.Range("B5").Value = 7
$ws.Range("B8").Value = 5
$ws.Range("B11").Value = 4
$ws.Range("B14").Value = 7
$ws.Range("B17").Value = 7
$ws.Range("B20").Value = 7
$ws.Range("B23").Value = 7
$ws.Range("B26").Value = 5
$ws.Range("B29").Value = 27

# Update the selected cell on the sheet
$ws.Activate()
$ws.Range("F32").Select()
